$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 761
$ws.Range("I33").Value = 654.1177
$ws.Range("J33").Value = 1366.6666
$ws.Range("K33").Value = 654.1177
$ws.Range("L33").Value = 1366.6666
$ws.Range("M33").Value = -425.1177
$ws.Range("N33").Value = -1824.6666
$ws.Range("H43").Value = 7328.125
$ws.Range("I43").Value = 1760.25
$ws.Range("J43").Value = 9184.083000000001
$ws.Range("K43").Value = 1760.25
$ws.Range("L43").Value = 9184.083000000001
$ws.Range("M43").Value = -1691.25
$ws.Range("N43").Value = -9322.083000000001
$ws.Range("H76").Value = 26325340
$ws.Range("I76").Value = 50013308
$ws.Range("K76").Value = 50013308
$ws.Range("M76").Value = -50012993
$ws.Range("H79").Value = 26325340
$ws.Range("I79").Value = 50013308
$ws.Range("K79").Value = 50013308
$ws.Range("M79").Value = -50012216
$ws.Range("H125").Value = 1752.8948
$ws.Range("I125").Value = 913.1429000000001
$ws.Range("J125").Value = 2242.75
$ws.Range("K125").Value = 8218.286100000001
$ws.Range("L125").Value = 20184.75
$ws.Range("M125").Value = -5758.286100000001
$ws.Range("N125").Value = -25104.75
$ws.Range("H127").Value = 671.4286
$ws.Range("I127").Value = 504.4
$ws.Range("J127").Value = 1089
$ws.Range("K127").Value = 1513.2
$ws.Range("L127").Value = 3267
$ws.Range("M127").Value = 3446.8
$ws.Range("N127").Value = -13187
$ws.Range("H131").Value = 1388.6111
$ws.Range("I131").Value = 466.1111
$ws.Range("J131").Value = 2311.111
$ws.Range("K131").Value = 1398.3333
$ws.Range("L131").Value = 6933.333
$ws.Range("M131").Value = 3641.6667
$ws.Range("N131").Value = -17013.333
$ws.Range("H141").Value = 4750.25
$ws.Range("I141").Value = 2000.5
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 6001.5
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = -821.5
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2965.9092
$ws.Range("I63").Value = 2965.9092
$ws.Range("K63").Value = 2965.9092
$ws.Range("M63").Value = -2279.9092
$ws.Range("H66").Value = 2965.9092
$ws.Range("I66").Value = 2965.9092
$ws.Range("K66").Value = 14829.546
$ws.Range("M66").Value = -11397.546
$ws.Range("H88").Value = 2103.5334
$ws.Range("I88").Value = 1982.091
$ws.Range("K88").Value = 1982.091
$ws.Range("M88").Value = -1576.091
$ws.Range("H91").Value = 2103.5334
$ws.Range("I91").Value = 1982.091
$ws.Range("K91").Value = 1982.091
$ws.Range("M91").Value = -578.0909999999999
$ws.Range("H132").Value = 30939.73
$ws.Range("I132").Value = 76335.36
$ws.Range("J132").Value = 3307.6086
$ws.Range("K132").Value = 229006.08
$ws.Range("L132").Value = 9922.825800000001
$ws.Range("M132").Value = -226476.08
$ws.Range("N132").Value = -14982.8258

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 34390
$ws.Range("J132").Value = 34390
$ws.Range("L132").Value = 34390
$ws.Range("N132").Value = -44510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 48006.273
$ws.Range("I99").Value = 144656
$ws.Range("J99").Value = 2903.0667
$ws.Range("K99").Value = 144656
$ws.Range("L99").Value = 2903.0667
$ws.Range("M99").Value = -143158
$ws.Range("N99").Value = -5899.066699999999
$ws.Range("H122").Value = 2258.375
$ws.Range("I122").Value = 2802
$ws.Range("J122").Value = 1352.3334
$ws.Range("K122").Value = 8406
$ws.Range("L122").Value = 4057.0002
$ws.Range("M122").Value = -5956
$ws.Range("N122").Value = -8957.0002
$ws.Range("H126").Value = 48006.273
$ws.Range("I126").Value = 144656
$ws.Range("J126").Value = 2903.0667
$ws.Range("K126").Value = 433968
$ws.Range("L126").Value = 8709.2001
$ws.Range("M126").Value = -431498
$ws.Range("N126").Value = -13649.2001
$ws.Range("H132").Value = 2628.818
$ws.Range("I132").Value = 1034.3334
$ws.Range("K132").Value = 3103.0002
$ws.Range("M132").Value = -573.0001999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 67.4375
$ws.Range("J2").Value = 90.90000000000001
$ws.Range("L2").Value = 545.4000000000001
$ws.Range("N2").Value = -771.4000000000001
$ws.Range("H32").Value = 676.0909
$ws.Range("I32").Value = 648.55554
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 1945.66662
$ws.Range("L32").Value = 2400
$ws.Range("M32").Value = -1662.66662
$ws.Range("N32").Value = -2966
$ws.Range("H38").Value = 105.92308
$ws.Range("I38").Value = 70.625
$ws.Range("K38").Value = 211.875
$ws.Range("M38").Value = 135.125
$ws.Range("H113").Value = 23810104
$ws.Range("I113").Value = 27778352
$ws.Range("J113").Value = 617.6667
$ws.Range("K113").Value = 83335056
$ws.Range("L113").Value = 1853.0001
$ws.Range("M113").Value = -83332886
$ws.Range("N113").Value = -6193.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4515.385
$ws.Range("I80").Value = 5075.8623
$ws.Range("J80").Value = 2890
$ws.Range("K80").Value = 5075.8623
$ws.Range("L80").Value = 2890
$ws.Range("M80").Value = -4077.8623
$ws.Range("N80").Value = -4886
$ws.Range("H83").Value = 4515.385
$ws.Range("I83").Value = 5075.8623
$ws.Range("J83").Value = 2890
$ws.Range("K83").Value = 25379.3115
$ws.Range("L83").Value = 14450
$ws.Range("M83").Value = -20387.3115
$ws.Range("N83").Value = -24434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5185.3447
$ws.Range("I40").Value = 7967.1875
$ws.Range("J40").Value = 1761.5385
$ws.Range("K40").Value = 7967.1875
$ws.Range("L40").Value = 1761.5385
$ws.Range("M40").Value = -7831.1875
$ws.Range("N40").Value = -2033.5385
$ws.Range("H122").Value = 3226.923
$ws.Range("I122").Value = 2908.3333
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8724.999899999999
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6274.999899999999
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 21744000
$ws.Range("I132").Value = 55560476
$ws.Range("K132").Value = 166681428
$ws.Range("M132").Value = -166678898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 385
$ws.Range("I113").Value = 380.4
$ws.Range("J113").Value = 396.5
$ws.Range("K113").Value = 1141.2
$ws.Range("L113").Value = 1189.5
$ws.Range("M113").Value = 1028.8
$ws.Range("N113").Value = -5529.5
$ws.Range("H122").Value = 86009.836
$ws.Range("I122").Value = 113701.445
$ws.Range("J122").Value = 2935
$ws.Range("K122").Value = 341104.335
$ws.Range("L122").Value = 8805
$ws.Range("M122").Value = -338654.335
$ws.Range("N122").Value = -13705
